$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @(3, 44334, 50, 26000, 28000, 27200, 2092),
    @(4, 44453, 50, 14000, 15000, 14600, 1123),
    @(5, 44462, 60, 14000, 15000, 14500, 1115),
    @(6, 44503, 35, 15000, 16000, 15429, 1187),
    @(7, 44377, 40, 14000, 15000, 14500, 1115),
    @(8, 44362, 40, 15000, 16000, 15500, 1192),
    @(9, 44435, 100, 13000, 14000, 13500, 1038),
    @(10, 44509, 100, 15000, 16000, 15500, 1192),
    @(11, 44159, 60, 30000, 32000, 31000, 2385),
    @(12, 44523, 40, 15000, 16000, 15500, 1192),
    @(13, 44355, 60, 18000, 20000, 19000, 1462),
    @(14, 44320, 50, 26000, 28000, 26800, 2062),
    @(15, 44488, 40, 16000, 17000, 16500, 1269),
    @(16, 44383, 50, 15000, 16000, 15400, 1185),
    @(17, 44327, 50, 24000, 25000, 24400, 1877),
    @(18, 44316, 50, 27000, 28000, 27400, 2108),
    @(19, 44425, 60, 14000, 15000, 14500, 1115),
    @(20, 44467, 100, 13000, 14000, 13500, 1038),
    @(21, 44264, 40, 30000, 32000, 31000, 2385),
    @(22, 44510, 40, 15000, 16000, 15500, 1192),
    @(23, 44308, 50, 26000, 27000, 26400, 2031),
    @(24, 44474, 40, 13000, 14000, 13500, 1038),
    @(25, 44350, 40, 23000, 25000, 24000, 1846),
    @(26, 44313, 50, 25000, 26000, 25600, 1969)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]   # D - Fecha
    $ws.Cells.Item($r, 10).Value = $row[2]  # J - Volumen
    $ws.Cells.Item($r, 11).Value = $row[3]  # K - Precio minimo
    $ws.Cells.Item($r, 12).Value = $row[4]  # L - Precio maximo
    $ws.Cells.Item($r, 13).Value = $row[5]  # M - Precio promedio ponderado
    $ws.Cells.Item($r, 16).Value = $row[6]  # P - Precio $/Kg
}

Write-Output "Updated rows 3-26 with permuted Jengibre price data"
